$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for the three new "大量成交區" columns ---
$ws.Range("G1").Value = "大量成交區"
$ws.Range("H1").Value = "大量成交區_2"
$ws.Range("I1").Value = "大量成交區_3"

# --- Update existing data rows 2-8 (column B-F) with the 2021/11/11 refreshed numbers ---
# Row 2 (Stock 1711 stays, only Previous_N_High/Low change)
$ws.Range("E2").Value = 29.35
$ws.Range("F2").Value = 27.5

# Row 3
$ws.Range("B3").Value = 3035
$ws.Range("C3").Value = 121.5
$ws.Range("D3").Value = 86.7
$ws.Range("E3").Value = 176
$ws.Range("F3").Value = 163

# Row 4
$ws.Range("B4").Value = 3141
$ws.Range("C4").Value = 185
$ws.Range("D4").Value = 127
$ws.Range("E4").Value = 220
$ws.Range("F4").Value = 216

# Row 5
$ws.Range("B5").Value = 3189
$ws.Range("C5").Value = 229
$ws.Range("D5").Value = 191
$ws.Range("E5").Value = 248
$ws.Range("F5").Value = 233

# Row 6
$ws.Range("B6").Value = 3504
$ws.Range("C6").Value = 121
$ws.Range("D6").Value = 99.9
$ws.Range("E6").Value = 137
$ws.Range("F6").Value = 130

# Row 7
$ws.Range("B7").Value = 6170
$ws.Range("C7").Value = 47.15
$ws.Range("D7").Value = 33.9
$ws.Range("E7").Value = 51.7
$ws.Range("F7").Value = 51.4

# Row 8
$ws.Range("B8").Value = 6411
$ws.Range("C8").Value = 227
$ws.Range("D8").Value = 190.5
$ws.Range("E8").Value = 227
$ws.Range("F8").Value = 204

# --- New row 9 (Stock 3033, newly tracked position) ---
$ws.Range("B9").Value = 3033
$ws.Range("C9").Value = 29.4
$ws.Range("D9").Value = 25.7
$ws.Range("E9").Value = 29.4
$ws.Range("F9").Value = 25.7

# --- Column widths for the newly added columns G:I (bestFit-equivalent, ~13.125 / 15.5 chars) ---
$ws.Columns.Item(7).ColumnWidth = 12.428571428571429
$ws.Columns.Item(8).ColumnWidth = 14.785714285714286
$ws.Columns.Item(9).ColumnWidth = 14.785714285714286

# --- Selection moves to F9, matching the last edited cell ---
$ws.Range("F9").Select()
